# Auto-generated edits applying the KPI re-run results diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 8150.0704200518003
$ws.Range("I2").Value = 7400.5411534089699
$ws.Range("J2").Value = 6473.3609008049798
$ws.Range("AB2").Value = 8316.3983878079598
$ws.Range("AC2").Value = 8390.6362283548497
$ws.Range("AD2").Value = 8183.7685218773504
$ws.Range("AM2").Value = 1.1250977419955099
$ws.Range("AN2").Value = 32.6586752295228
$ws.Range("AR2").Value = 56023.972474265604
$ws.Range("AS2").Value = 32333.204041784898
$ws.Range("AT2").Value = 3902.9800741229501
$ws.Range("F3").Value = 11374.1803313611
$ws.Range("G3").Value = 14631.2853286749
$ws.Range("H3").Value = 12997.503750198601
$ws.Range("I3").Value = 11374.3598271231
$ws.Range("J3").Value = 8892.5799839159008
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = 0.91500000000000004
$ws.Range("S3").Value = 0.755
$ws.Range("T3").Value = 0.63200000000000001
$ws.Range("Z3").Value = 11998.0805183134
$ws.Range("AA3").Value = 14631.2853286749
$ws.Range("AB3").Value = 14204.922131364599
$ws.Range("AC3").Value = 15065.377254467699
$ws.Range("AD3").Value = 14070.537912056199
$ws.Range("AM3").Value = [double]"7.70167958486163E-13"
$ws.Range("AN3").Value = 130.070431977254
$ws.Range("AR3").Value = 84582.909221272799
$ws.Range("AS3").Value = 57449.8073820848
$ws.Range("AT3").Value = 460.646667415038
$ws.Range("G4").Value = 14632.6487718975
$ws.Range("H4").Value = 27436.348027577002
$ws.Range("I4").Value = 23783.9556423457
$ws.Range("J4").Value = 19251.780761506499
$ws.Range("S4").Value = 0.91700000000000004
$ws.Range("T4").Value = 0.80800000000000005
$ws.Range("AA4").Value = 23116.348770770201
$ws.Range("AB4").Value = 27436.348027577002
$ws.Range("AC4").Value = 25936.7019000498
$ws.Range("AD4").Value = 23826.461338498098
$ws.Range("AL4").Value = 0
$ws.Range("AM4").Value = 0
$ws.Range("AN4").Value = 23.992015566669501
$ws.Range("AR4").Value = 85104.733203326701
$ws.Range("AS4").Value = 57788.213950614001
$ws.Range("AT4").Value = 277.22925389091398
$ws.Range("H5").Value = 1221.6869502422001
$ws.Range("I5").Value = 38882.7580529392
$ws.Range("J5").Value = 45193.622994310099
$ws.Range("T5").Value = 0.998
$ws.Range("AB5").Value = 13726.819665642701
$ws.Range("AC5").Value = 42217.978341953502
$ws.Range("AD5").Value = 45284.1913770643
$ws.Range("AN5").Value = [double]"1.1757081824619299E-6"
$ws.Range("AO5").Value = [double]"-2.8950539600761901E-12"
$ws.Range("AR5").Value = 85298.067997491802
$ws.Range("AS5").Value = 57931.303305896203
$ws.Range("AT5").Value = 226.98381500830499
$ws.Range("I6").Value = 1781.60096455194
$ws.Range("J6").Value = 83811.460958038806
$ws.Range("AC6").Value = 14028.3540515901
$ws.Range("AD6").Value = 83811.460958038806
$ws.Range("AO6").Value = 0
$ws.Range("AR6").Value = 85593.061922590699
$ws.Range("AS6").Value = 58224.340284935402
$ws.Range("AT6").Value = 225.02686894876501
$ws.Range("I7").Value = 8441.6118179445803
$ws.Range("J7").Value = 7669.3089259242397
$ws.Range("K7").Value = 6033.2378625561596
$ws.Range("L7").Value = 5285.2230851219301
$ws.Range("V7").Value = 0.94199999999999995
$ws.Range("X7").Value = 0.67900000000000005
$ws.Range("AG7").Value = 8441.6118179445803
$ws.Range("AH7").Value = 8141.5169062053601
$ws.Range("AI7").Value = 7295.3299426313897
$ws.Range("AJ7").Value = 7783.8337041560399
$ws.Range("AU7").Value = [double]"-2.1827872842550201E-13"
$ws.Range("AV7").Value = 32.776565429385698
$ws.Range("AZ7").Value = 69929.381691546805
$ws.Range("BA7").Value = 48164.348875568998
$ws.Range("BB7").Value = 2093.4516859713899
$ws.Range("F8").Value = 4378.8272603183304
$ws.Range("G8").Value = 17000.0000001515
$ws.Range("H8").Value = 17000.0000000418
$ws.Range("I8").Value = 16252.000000067001
$ws.Range("J8").Value = 12867.894606718701
$ws.Range("K8").Value = 10046.869849397801
$ws.Range("L8").Value = 7510.91411516469
$ws.Range("U8").Value = 0.95599999999999996
$ws.Range("W8").Value = 0.72799999999999998
$ws.Range("X8").Value = 0.61299999999999999
$ws.Range("AD8").Value = 4378.8272603183304
$ws.Range("AE8").Value = 17000.0000001515
$ws.Range("AF8").Value = 17000.0000000418
$ws.Range("AG8").Value = 17000.000000004798
$ws.Range("AH8").Value = 14063.2727942281
$ws.Range("AI8").Value = 13800.6453975245
$ws.Range("AJ8").Value = 12252.714706630801
$ws.Range("AT8").Value = 0
$ws.Range("AU8").Value = [double]"-7.8676779678870597E-10"
$ws.Range("AV8").Value = 56.092483562184697
$ws.Range("AZ8").Value = 85056.505831859904
$ws.Range("BA8").Value = 61578.074915619698
$ws.Range("BB8").Value = 380.05358570879702
$ws.Range("H9").Value = 426.55068922957099
$ws.Range("I9").Value = 6132.5855215791198
$ws.Range("J9").Value = 30492.7969034815
$ws.Range("K9").Value = 25191.975294191299
$ws.Range("L9").Value = 18499.746147696002
$ws.Range("W9").Value = 0.93200000000000005
$ws.Range("X9").Value = 0.80500000000000005
$ws.Range("AF9").Value = 21327.534461478499
$ws.Range("AG9").Value = 6151.0386374915997
$ws.Range("AH9").Value = 30553.904712907301
$ws.Range("AI9").Value = 27030.016410076401
$ws.Range("AJ9").Value = 22981.051113601399
$ws.Range("AU9").Value = [double]"5.7919039222264902E-14"
$ws.Range("AV9").Value = 122.515584251492
$ws.Range("AW9").Value = 0
$ws.Range("AZ9").Value = 80743.654556177498
$ws.Range("BA9").Value = 57316.646779459297
$ws.Range("BB9").Value = 431.47672523089102
$ws.Range("J10").Value = 1351.9792331163201
$ws.Range("K10").Value = 36646.350538423503
$ws.Range("L10").Value = 42195.932829158497
$ws.Range("AH10").Value = 27591.412920741201
$ws.Range("AI10").Value = 40050.656326145901
$ws.Range("AJ10").Value = 43590.839699543903
$ws.Range("AU10").Value = 0
$ws.Range("AV10").Value = 1.49973408624987
$ws.Range("AZ10").Value = 80194.262600697897
$ws.Range("BA10").Value = 56734.522389030899
$ws.Range("BB10").Value = 398.74429028172699
$ws.Range("K11").Value = 2677.0875948598
$ws.Range("L11").Value = 77551.778532794197
$ws.Range("AI11").Value = 32647.409693412199
$ws.Range("AJ11").Value = 77551.778532794197
$ws.Range("AV11").Value = 0
$ws.Range("AW11").Value = 0
$ws.Range("AX11").Value = 0
$ws.Range("AZ11").Value = 80228.866127653993
$ws.Range("BA11").Value = 56754.024200859902
$ws.Range("BB11").Value = 383.642575155072

# Cells that now need scientific-notation number format (style index 1)
$ws.Range("AN5").NumberFormat = "0.00E+00"
$ws.Range("AO5").NumberFormat = "0.00E+00"
$ws.Range("AU7").NumberFormat = "0.00E+00"

# Cells that revert to default/General formatting (clear the style)
$ws.Range("AO6").ClearFormats()
$ws.Range("AT8").ClearFormats()
$ws.Range("AW9").ClearFormats()
$ws.Range("AU10").ClearFormats()
$ws.Range("AV11").ClearFormats()
$ws.Range("AW11").ClearFormats()
$ws.Range("AX11").ClearFormats()

# Update the active selection to match the saved view
$ws.Range("D18").Select()
